# WeeklyTracker update by Khaleel
# - Fix typo "In Progrees " -> "In Progress " on the Khaleel sheet (Task 1 status)
# - Fill in the "Task 2" row (pendios dashboard work) details on the Khaleel sheet
# - Widen column D on the Khaleel sheet to fit the new long text
# - Make the Khaleel sheet the active tab, with D8 selected

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Khaleel")

# Fix the "In Progrees " typo in the Task 1 status cell (G6)
$ws.Range("G6").Value = "In Progress "

# Fill in Task 2 details (row 7): Task Name, Complexity level, Team Member, Status, Remarks
$ws.Range("D7").Value = "Currently working on pendios dashboard functionalities with Arun and Rahamath"
$ws.Range("E7").Value = "High"
$ws.Range("F7").Value = "Arun, Rahamath and Khaleel"
$ws.Range("G7").Value = "In Progress"
$ws.Range("H7").Value = "Action page is under progress"

# Widen column D so the new long task name is readable
$ws.Columns.Item(4).ColumnWidth = 68.5

# Make Khaleel the active sheet/tab and select D8
$null = $ws.Activate()
$ws.Range("D8").Select() | Out-Null
